# update scripts wuth new tpm
# Refresh the NATMI Fn1-Sdc2 ligand-receptor stats with the recomputed TPM
# values (ligand/receptor average & total expression for the "ECs" cluster
# changed, which ripples into every derived specificity / edge-weight column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 6.968297403059999
$ws.Range("R2").Value = 62.71467662753999
$ws.Range("S2").Value = 0.0001415962350341598
$ws.Range("T2").Value = 0.0001415962350341598
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 698.4678102567199
$ws.Range("R3").Value = 6286.210292310479
$ws.Range("S3").Value = 0.01419290918632077
$ws.Range("T3").Value = 0.01419290918632077
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 147.03694402614
$ws.Range("R4").Value = 1323.33249623526
$ws.Range("S4").Value = 0.002987799814038828
$ws.Range("T4").Value = 0.002987799814038829
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 385.892470936405
$ws.Range("R5").Value = 3473.032238427645
$ws.Range("S5").Value = 0.007841358921998552
$ws.Range("T5").Value = 0.007841358921998555
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("S6").Value = 0.7859791967662327
$ws.Range("T6").Value = 0.7859791967662328
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("S7").Value = 0.1654592773833777
$ws.Range("T7").Value = 0.1654592773833777
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 1.116695
$ws.Range("N8").Value = 3.350085
$ws.Range("O8").Value = 0.008174214292497491
$ws.Range("P8").Value = 0.008174214292497492
$ws.Range("Q8").Value = 9.412330325373333
$ws.Range("R8").Value = 84.71097292835999
$ws.Range("S8").Value = 0.0001912591354647778
$ws.Range("T8").Value = 0.0001912591354647779
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.8193429796700005
$ws.Range("P9").Value = 0.8193429796700005
$ws.Range("Q9").Value = 943.4456326289244
$ws.Range("S9").Value = 0.01917087371744694
$ws.Range("T9").Value = 0.01917087371744695
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.172482806037502
$ws.Range("P10").Value = 0.1724828060375021
$ws.Range("Q10").Value = 198.6080970940933
$ws.Range("S10").Value = 0.004035728840085548
$ws.Range("T10").Value = 0.004035728840085549
